$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "Carol Poole" above the current row 18
# (Camille Fritsch / row 17, Caroline Lötter / row 18 -> 19), keeping the
# contributor list alphabetically ordered.
$ws.Rows("18:18").Insert() | Out-Null
$ws.Range("A18").Value = "Carol Poole"
$ws.Range("B18").Value = "South African National Biodiversity Institute"

# The data range grew by one row (A1:B167 -> A1:B168): refresh the
# AutoFilter so its button range covers the new last row too. (Existing
# row order must NOT be touched, so this intentionally does not re-run a
# Sort over the data -- only the filter's own range is refreshed.)
$ws.AutoFilterMode = $false
$ws.Range("A1:B168").AutoFilter() | Out-Null

# Keep the hidden _xlnm._FilterDatabase name in sync with the AutoFilter.
$filterName = $wb.Names.Item(1)
$filterName.RefersTo = "=Sheet1!`$A`$1:`$B`$168"

# Restore the selection to the cell the author last touched.
$ws.Range("B24").Select() | Out-Null
